# Generate Report for Archive
#
# 1) Update status text from "Ready for handoff" to "In Translation"
#    everywhere it appears (Overview!E2/F2, zh-cn!C2, de-de!C2).
# 2) Narrow the "Status" columns (Overview E & F, zh-cn C, de-de C)
#    from their current width down to the newly computed narrower width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# --- 1. Text update ------------------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# --- 2. Column width update ----------------------------------------------
# Original stored width 17.2159881591797 -> new stored width 13.4101845877511
# (columns addressed numerically: Overview E=5, F=6 ; zh-cn/de-de C=3)
$newColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth     = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth     = $newColumnWidth

Write-Output "edit complete"
